$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.08"
$ws.Range("D3").Value = "'21.66"
$ws.Range("D4").Value = "'5.309"
$ws.Range("D5").Value = "'0.05596"
$ws.Range("D6").Value = "'3.377"
$ws.Range("D7").Value = "'6.374"
$ws.Range("D8").Value = "'0.8160"
$ws.Range("D9").Value = "'0.9761"
$ws.Range("D10").Value = "'0.1406"
$ws.Range("D11").Value = "'0.07397"
$ws.Range("D12").Value = "'0.03119"
$ws.Range("D13").Value = "'0.03032"
$ws.Range("D14").Value = "'0.09305"
$ws.Range("D15").Value = "'3.568"
$ws.Range("D16").Value = "'0.001611"
$ws.Range("D17").Value = "'0.04675"
$ws.Range("D18").Value = "'0.0005747"
$ws.Range("D19").Value = "'0.006341"
$ws.Range("D20").Value = "'0.005058"
$ws.Range("D21").Value = "'0.001033"
$ws.Range("D22").Value = "'0.0001496"
$ws.Range("D23").Value = "'3.780"
$ws.Range("D24").Value = "'2.120"
$ws.Range("D25").Value = "'0.3250"
$ws.Range("D26").Value = "'0.1297"
$ws.Range("D28").Value = "'0.0003080"
$ws.Range("D40").Value = "'0.03915"
$ws.Range("D41").Value = "'0.007027"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1047"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003006"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.007748"
$ws.Range("D45").Value = "'0.00005769"
$ws.Range("D46").Value = "'0.00000000745"
$ws.Range("D47").Value = "'0.0005468"
$ws.Range("D48").Value = "'0.6757"
$ws.Range("D49").Value = "'0.1446"
$ws.Range("D50").Value = "'0.00002087"
$ws.Range("D51").Value = "'0.01004"
